# Fix: the serial number recorded for the ADCPA-m (glider deployment
# CE05MOAS-GL383_00001) was wrong in the calibration-coefficient sheet.
# Correct the "Sensor Serial Number" column (D) for that instrument's
# four calibration-coefficient rows (9-12) on the Asset_Cal_Info sheet.

$wb = $excel.ActiveWorkbook

$wsAsset = $wb.Worksheets.Item("Asset_Cal_Info")

$wsAsset.Range("D9").Value = 649982
$wsAsset.Range("D10").Value = 649982
$wsAsset.Range("D11").Value = 649982
$wsAsset.Range("D12").Value = 649982

# Reflect this was the last-edited/active sheet, with the corrected
# rows selected, before the file was saved.
$wsAsset.Activate()
$wsAsset.Range("D9:D12").Select()

$wb.Save()
